# The observation records in rows 2-18 got re-sorted: every record kept its
# own field values, but which spreadsheet row it landed on changed (a
# permutation of rows 2..18). Reproduce that by copying each record's data
# into its new row.
#
# Only columns A, B, D, E, F, G, H, M, Q and R ever actually differ between
# any two rows in this sheet (every other column - C, I, K, L, N, P, S, T,
# U, V, W, Y, Z, AA, AB, ... - already holds the same value/blank in every
# row, so the shuffle never changes what they display and they're left
# untouched). Restricting the writes to exactly those columns also avoids
# Excel's autoconversion of the untouched "2023-08-08" date-as-text cells
# into real dates, which a blanket whole-row Value2 copy would trigger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 18
$rowCount = $lastRow - $firstRow + 1

# Destination row offset (0-based, row2 = offset 0) -> source row offset it
# should copy its data from, e.g. new row 2 gets what used to be in row 14.
$rowMap = @{
    0  = 12   # row 2  <- row 14
    1  = 13   # row 3  <- row 15
    2  = 4    # row 4  <- row 6
    3  = 0    # row 5  <- row 2
    4  = 8    # row 6  <- row 10
    5  = 2    # row 7  <- row 4
    6  = 3    # row 8  <- row 5
    7  = 5    # row 9  <- row 7
    8  = 1    # row 10 <- row 3
    9  = 14   # row 11 <- row 16
    10 = 7    # row 12 <- row 9
    11 = 10   # row 13 <- row 12
    12 = 6    # row 14 <- row 8
    13 = 9    # row 15 <- row 11
    14 = 11   # row 16 <- row 13
    15 = 16   # row 17 <- row 18
    16 = 15   # row 18 <- row 17
}

function Permute-ColumnBlock($firstCol, $lastCol) {
    $range = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
    $colCount = $lastCol - $firstCol + 1
    $src = $range.Value2

    $dst = New-Object 'object[,]' $rowCount, $colCount
    for ($r = 0; $r -lt $rowCount; $r++) {
        $sourceOffset = $rowMap[$r]
        for ($c = 0; $c -lt $colCount; $c++) {
            $dst[$r, $c] = $src[$sourceOffset + 1, $c + 1]
        }
    }

    $range.Value2 = $dst
}

Permute-ColumnBlock 1  2    # A:B
Permute-ColumnBlock 4  8    # D:H
Permute-ColumnBlock 13 13   # M
Permute-ColumnBlock 17 18   # Q:R
